# Update cryptocurrency price (D) and volume-change (E) columns with freshly
# scraped figures. Price-column values must be written as literal text so
# that numeric-looking strings (e.g. "1.001", "102.00") keep their exact
# displayed digits instead of being re-interpreted as floating point numbers
# (which would drop trailing zeros / introduce binary rounding noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    "D2" = "30.335.92"
    "D3" = "1.938.09"
    "D4" = "1.001"
    "D5" = "251.45"
    "D6" = "0.7235"
    "D9" = "28.12"
    "D10" = "0.07230"
    "D11" = "0.8116"
    "D12" = "0.08105"
    "D13" = "1.938.12"
    "D14" = "5.487"
    "D15" = "94.66"
    "D16" = "15.15"
    "D17" = "30.348.43"
    "D18" = "0.000008294"
    "D19" = "249.79"
    "D21" = "2.192.12"
    "D23" = "1.002"
    "D24" = "6.995"
    "D26" = "163.64"
    "D27" = "2.385"
    "D28" = "19.33"
    "D29" = "0.1321"
    "D30" = "1.570"
    "D31" = "1.346"
    "D32" = "4.443"
    "D33" = "4.182"
    "D34" = "0.05208"
    "D35" = "1.290"
    "D36" = "0.7514"
    "D37" = "2.753"
    "D38" = "0.01982"
    "D39" = "2.834"
    "D40" = "80.75"
    "D42" = "0.4547"
    "D43" = "2.039"
    "D46" = "102.00"
    "D47" = "9.816"
    "D48" = "7.461"
    "D49" = "36.87"
    "D51" = "0.06046"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like "1.001" / "102.00" keep their exact
    # textual form instead of being coerced into a Double.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    # Restore the default style so no stray format is left on the cell.
    $cell.Style = "Normal"
}

$volumeUpdates = @{
    "E3" = "  -2.94%  "
    "E4" = "  +0.06%  "
    "E5" = "  -1.28%  "
    "E6" = "  -6.31%  "
    "E7" = "  +0.09%  "
    "E8" = "  -4.64%  "
    "E9" = "  +0.83%  "
    "E10" = "  +2.18%  "
    "E11" = "  -3.62%  "
    "E12" = "  -1.23%  "
    "E13" = "  -2.93%  "
    "E14" = "  -2.62%  "
    "E15" = "  -6.06%  "
    "E16" = "  -0.51%  "
    "E18" = "  +3.51%  "
    "E19" = "  -8.32%  "
    "E20" = "  -1.52%  "
    "E21" = "  -2.78%  "
    "E22" = "  +0.08%  "
    "E23" = "  +0.04%  "
    "E24" = "  -1.43%  "
    "E25" = "  -2.40%  "
    "E26" = "  -1.25%  "
    "E27" = "  -0.32%  "
    "E28" = "  -2.93%  "
    "E29" = "  -7.05%  "
    "E30" = "  -1.65%  "
    "E31" = "  -1.74%  "
    "E32" = "  -3.54%  "
    "E33" = "  -5.85%  "
    "E34" = "  -1.88%  "
    "E35" = "  +3.92%  "
    "E36" = "  -4.93%  "
    "E37" = "  -0.57%  "
    "E38" = "  -1.05%  "
    "E39" = "  -2.94%  "
    "E40" = "  -1.64%  "
    "E41" = "  -4.86%  "
    "E42" = "  -2.77%  "
    "E43" = "  -4.56%  "
    "E44" = "  -0.88%  "
    "E46" = "  -3.23%  "
    "E47" = "  -2.34%  "
    "E48" = "  -3.42%  "
    "E49" = "  -1.67%  "
    "E50" = "  -3.07%  "
    "E51" = "  +0.59%  "
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
